$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 6
$ws.Range("H6").Value = 419.9
$ws.Range("I6").Value = 139.8
$ws.Range("J6").Value = 700
$ws.Range("K6").Value = 419.4
$ws.Range("L6").Value = 2100
$ws.Range("M6").Value = -307.4
$ws.Range("N6").Value = -2324

# Row 129
$ws.Range("H129").Value = 978.6857
$ws.Range("I129").Value = 454.5
$ws.Range("J129").Value = 1252.174
$ws.Range("K129").Value = 1363.5
$ws.Range("L129").Value = 3756.522
$ws.Range("M129").Value = 3636.5
$ws.Range("N129").Value = -13756.522

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 28188.486
$ws.Range("I32").Value = 11087.635
$ws.Range("J32").Value = 182096.14
$ws.Range("K32").Value = 11087.635
$ws.Range("L32").Value = 182096.14
$ws.Range("M32").Value = -10800.635
$ws.Range("N32").Value = -182670.14

# Row 74
$ws.Range("H74").Value = 23127.55
$ws.Range("I74").Value = 2404.1875
$ws.Range("K74").Value = 2404.1875
$ws.Range("M74").Value = -1530.1875

# Row 77
$ws.Range("H77").Value = 23127.55
$ws.Range("I77").Value = 2404.1875
$ws.Range("K77").Value = 12020.9375
$ws.Range("M77").Value = -7652.9375

# Row 97
$ws.Range("H97").Value = 1157.3704
$ws.Range("I97").Value = 697.7778
$ws.Range("J97").Value = 2076.5557
$ws.Range("K97").Value = 697.7778
$ws.Range("L97").Value = 2076.5557
$ws.Range("M97").Value = -201.7778
$ws.Range("N97").Value = -3068.5557

# Row 107
$ws.Range("H107").Value = 24076
$ws.Range("J107").Value = 24076
$ws.Range("L107").Value = 24076
$ws.Range("N107").Value = -31756

# Row 132
$ws.Range("H132").Value = 3430.1936
$ws.Range("I132").Value = 2063.5
$ws.Range("J132").Value = 8116
$ws.Range("K132").Value = 6190.5
$ws.Range("L132").Value = 24348
$ws.Range("M132").Value = -3660.5
$ws.Range("N132").Value = -29408

$ws = $wb.Worksheets.Item("BSM")
# Row 7
$ws.Range("H7").Value = 100
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()

# Row 94
$ws.Range("H94").Value = 637
$ws.Range("I94").Value = 538.5833
$ws.Range("J94").Value = 805.7143
$ws.Range("K94").Value = 538.5833
$ws.Range("L94").Value = 805.7143
$ws.Range("M94").Value = -87.58330000000001
$ws.Range("N94").Value = -1707.7143

# Row 134
$ws.Range("H134").Value = 2586.6875
$ws.Range("I134").Value = 1997.3
$ws.Range("J134").Value = 3569
$ws.Range("K134").Value = 5991.9
$ws.Range("L134").Value = 10707
$ws.Range("M134").Value = -3456.9
$ws.Range("N134").Value = -15777

$ws = $wb.Worksheets.Item("CRP")
# Row 2
$ws.Range("H2").Value = 1871.75
$ws.Range("I2").Value = 987
$ws.Range("J2").Value = 2166.6667
$ws.Range("K2").Value = 987
$ws.Range("L2").Value = 2166.6667
$ws.Range("M2").Value = -874
$ws.Range("N2").Value = -2392.6667

# Row 16
$ws.Range("H16").Value = 2943217.5
$ws.Range("I16").Value = 7353891
$ws.Range("J16").Value = 2768.8333
$ws.Range("K16").Value = 7353891
$ws.Range("L16").Value = 2768.8333
$ws.Range("M16").Value = -7353604
$ws.Range("N16").Value = -3342.8333

# Row 31
$ws.Range("H31").Value = 3319.2334
$ws.Range("I31").Value = 1660.5853
$ws.Range("J31").Value = 6898.421
$ws.Range("K31").Value = 1660.5853
$ws.Range("L31").Value = 6898.421
$ws.Range("M31").Value = -1365.5853
$ws.Range("N31").Value = -7488.421

# Row 34
$ws.Range("H34").Value = 3319.2334
$ws.Range("I34").Value = 1660.5853
$ws.Range("J34").Value = 6898.421
$ws.Range("K34").Value = 1660.5853
$ws.Range("L34").Value = 6898.421
$ws.Range("M34").Value = -1458.5853
$ws.Range("N34").Value = -7302.421

# Row 113
$ws.Range("H113").Value = 2943217.5
$ws.Range("I113").Value = 7353891
$ws.Range("J113").Value = 2768.8333
$ws.Range("K113").Value = 7353891
$ws.Range("L113").Value = 2768.8333
$ws.Range("M113").Value = -7351721
$ws.Range("N113").Value = -7108.8333

# Row 131
$ws.Range("H131").Value = 47740
$ws.Range("J131").Value = 47740
$ws.Range("L131").Value = 47740
$ws.Range("N131").Value = -57820

$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 99.17778
$ws.Range("I4").Value = 97
$ws.Range("J4").Value = 113.333336
$ws.Range("K4").Value = 291
$ws.Range("L4").Value = 340.000008
$ws.Range("M4").Value = -179
$ws.Range("N4").Value = -564.000008

# Row 25
$ws.Range("H25").Value = 2800
$ws.Range("J25").Value = 5000
$ws.Range("L25").Value = 15000
$ws.Range("N25").Value = -15338

# Row 30
$ws.Range("H30").Value = 2800
$ws.Range("J30").Value = 5000
$ws.Range("L30").Value = 15000
$ws.Range("N30").Value = -15204

# Row 34
$ws.Range("H34").Value = 810.6667
$ws.Range("I34").Value = 274.7143
$ws.Range("J34").Value = 1151.7273
$ws.Range("K34").Value = 824.1428999999999
$ws.Range("L34").Value = 3455.1819
$ws.Range("M34").Value = -740.1428999999999
$ws.Range("N34").Value = -3623.1819

# Row 39
$ws.Range("H39").Value = 2900
$ws.Range("J39").Value = 2900
$ws.Range("L39").Value = 8700
$ws.Range("N39").Value = -9288

# Row 57
$ws.Range("H57").Value = 8813
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 8813
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 26439
$ws.Range("N57").Value = -27557
$ws.Range("M57").ClearContents()

# Row 131
$ws.Range("H131").Value = 897.98
$ws.Range("I131").Value = 405.55554
$ws.Range("J131").Value = 946.68134
$ws.Range("K131").Value = 1216.66662
$ws.Range("L131").Value = 2840.04402
$ws.Range("M131").Value = 3823.33338
$ws.Range("N131").Value = -12920.04402

# Row 137
$ws.Range("H137").Value = 3129.4614
$ws.Range("I137").Value = 2208.3333
$ws.Range("J137").Value = 3919
$ws.Range("K137").Value = 6624.999899999999
$ws.Range("L137").Value = 11757
$ws.Range("M137").Value = -1524.999899999999
$ws.Range("N137").Value = -21957

$ws = $wb.Worksheets.Item("GSM")
# Row 126
$ws.Range("H126").Value = 5409.1787
$ws.Range("I126").Value = 1756.9
$ws.Range("J126").Value = 7438.222
$ws.Range("K126").Value = 5270.700000000001
$ws.Range("L126").Value = 22314.666
$ws.Range("M126").Value = -2800.700000000001
$ws.Range("N126").Value = -27254.666

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 8098069.5
$ws.Range("I7").Value = 13494464
$ws.Range("J7").Value = 3477.8333
$ws.Range("K7").Value = 13494464
$ws.Range("L7").Value = 3477.8333
$ws.Range("M7").Value = -13494352
$ws.Range("N7").Value = -3701.8333

# Row 22
$ws.Range("H22").Value = 1129
$ws.Range("I22").Value = 891.3333
$ws.Range("J22").Value = 1366.6666
$ws.Range("K22").Value = 891.3333
$ws.Range("L22").Value = 1366.6666
$ws.Range("M22").Value = -596.3333
$ws.Range("N22").Value = -1956.6666

# Row 27
$ws.Range("H27").Value = 1129
$ws.Range("I27").Value = 891.3333
$ws.Range("J27").Value = 1366.6666
$ws.Range("K27").Value = 891.3333
$ws.Range("L27").Value = 1366.6666
$ws.Range("M27").Value = -784.3333
$ws.Range("N27").Value = -1580.6666

# Row 126
$ws.Range("H126").Value = 8098069.5
$ws.Range("I126").Value = 13494464
$ws.Range("J126").Value = 3477.8333
$ws.Range("K126").Value = 40483392
$ws.Range("L126").Value = 10433.4999
$ws.Range("M126").Value = -40480922
$ws.Range("N126").Value = -15373.4999

# Row 127
$ws.Range("H127").Value = 42658.75
$ws.Range("J127").Value = 42658.75
$ws.Range("L127").Value = 42658.75
$ws.Range("N127").Value = -52578.75

# Row 132
$ws.Range("H132").Value = 2941.32
$ws.Range("I132").Value = 2641.9
$ws.Range("J132").Value = 4139
$ws.Range("K132").Value = 7925.700000000001
$ws.Range("L132").Value = 12417
$ws.Range("M132").Value = -5395.700000000001
$ws.Range("N132").Value = -17477

$ws = $wb.Worksheets.Item("WVR")
# Row 54
$ws.Range("H54").Value = 10750
$ws.Range("J54").Value = 10750
$ws.Range("L54").Value = 10750
$ws.Range("N54").Value = -11790

# Row 62
$ws.Range("H62").Value = 25000
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 25000
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 25000
$ws.Range("N62").Value = -26248
$ws.Range("M62").ClearContents()

# Row 65
$ws.Range("H65").Value = 25000
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 25000
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 125000
$ws.Range("N65").Value = -131240
$ws.Range("M65").ClearContents()

# Row 81
$ws.Range("H81").Value = 11142.363
$ws.Range("I81").Value = 18761
$ws.Range("J81").Value = 2000
$ws.Range("K81").Value = 37522
$ws.Range("L81").Value = 4000
$ws.Range("M81").Value = -36461
$ws.Range("N81").Value = -6122

# Row 84
$ws.Range("H84").Value = 11142.363
$ws.Range("I84").Value = 18761
$ws.Range("J84").Value = 2000
$ws.Range("K84").Value = 187610
$ws.Range("L84").Value = 20000
$ws.Range("M84").Value = -182306
$ws.Range("N84").Value = -30608

# Row 100
$ws.Range("H100").Value = 1112032.5
$ws.Range("I100").Value = 516.6667
$ws.Range("J100").Value = 3335064.2
$ws.Range("K100").Value = 1033.3334
$ws.Range("L100").Value = 6670128.4
$ws.Range("M100").Value = -492.3334
$ws.Range("N100").Value = -6671210.4
